$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the ACE_landing_page_data table: shift the yearly series forward by
# one year and refresh all computed metrics for the new release (202312
# December full release).
$ws1.Range("A2").Value = [double]"2022"
$ws1.Range("B2").Value = [double]"470.15663229742336"
$ws1.Range("C2").Value = [double]"8921474338.786293"
$ws1.Range("D2").Value = [double]"18975536.504061323"
$ws1.Range("E2").Value = [double]"0.8836954720487562"
$ws1.Range("F2").Value = [double]"133.16679170159205"
$ws1.Range("G2").Value = [double]"319.46355316139619"
$ws1.Range("H2").Value = [double]"-0.34650568426651729"
$ws1.Range("I2").Value = [double]"3.4409572293988999E-2"
$ws1.Range("J2").Value = [double]"0.58288993092918751"
$ws1.Range("K2").Value = [double]"0.46856790462369879"
$ws1.Range("L2").Value = [double]"-2.3420464129504381E-2"
$ws1.Range("M2").Value = [double]"-0.35179038408891561"
$ws1.Range("N2").Value = [double]"96.927408574309311"
$ws1.Range("O2").Value = [double]"93.290786227371314"
$ws1.Range("A3").Value = [double]"2021"
$ws1.Range("B3").Value = [double]"719.45022470427932"
$ws1.Range("C3").Value = [double]"8624702030.7452507"
$ws1.Range("D3").Value = [double]"11987906.507764762"
$ws1.Range("E3").Value = [double]"0.60173960582040065"
$ws1.Range("F3").Value = [double]"136.36041593160243"
$ws1.Range("G3").Value = [double]"492.83988592544631"
$ws1.Range("H3").Value = [double]"-0.25294214079289401"
$ws1.Range("I3").Value = [double]"-4.9983832135574113E-2"
$ws1.Range("J3").Value = [double]"0.27167682684274408"
$ws1.Range("K3").Value = [double]"0.25122106103595776"
$ws1.Range("L3").Value = [double]"-8.2189450562490274E-2"
$ws1.Range("M3").Value = [double]"-0.246553986261564"
$ws1.Range("N3").Value = [double]"93.703124149707335"
$ws1.Range("O3").Value = [double]"58.937001496122846"
$ws1.Range("A4").Value = [double]"2020"
$ws1.Range("B4").Value = [double]"963.04485099436852"
$ws1.Range("C4").Value = [double]"9078479211.7096443"
$ws1.Range("D4").Value = [double]"9426849.8526687324"
$ws1.Range("E4").Value = [double]"0.48092189666483548"
$ws1.Range("F4").Value = [double]"148.57141924896419"
$ws1.Range("G4").Value = [double]"654.11439829654353"
$ws1.Range("H4").Value = [double]"1.2134441067149506"
$ws1.Range("I4").Value = [double]"-4.188692214530354E-2"
$ws1.Range("J4").Value = [double]"-0.56713924921435421"
$ws1.Range("K4").Value = [double]"-0.50702940127839669"
$ws1.Range("L4").Value = [double]"8.47184440380222E-2"
$ws1.Range("M4").Value = [double]"1.2196723111014927"
$ws1.Range("N4").Value = [double]"98.633188907032832"
$ws1.Range("O4").Value = [double]"46.345895633286567"
$ws1.Range("A5").Value = [double]"2019"
$ws1.Range("B5").Value = [double]"435.08885002913263"
$ws1.Range("C5").Value = [double]"9475373441.3449364"
$ws1.Range("D5").Value = [double]"21778019.456739664"
$ws1.Range("E5").Value = [double]"0.97555898447490952"
$ws1.Range("F5").Value = [double]"136.96772657048723"
$ws1.Range("G5").Value = [double]"294.68962379043467"
$ws1.Range("H5").Value = [double]"-2.4139029120445743E-3"
$ws1.Range("I5").Value = [double]"1.411651830421401E-2"
$ws1.Range("J5").Value = [double]"1.6570420602805447E-2"
$ws1.Range("K5").Value = [double]"1.0149067284991542E-2"
$ws1.Range("L5").Value = [double]"6.1435452387501588E-3"
$ws1.Range("M5").Value = [double]"-1.673078925594762E-3"
$ws1.Range("N5").Value = [double]"102.94524851688867"
$ws1.Range("O5").Value = [double]"107.06883345086935"
$ws1.Range("A6").Value = [double]"2018"
$ws1.Range("B6").Value = [double]"436.14165363691069"
$ws1.Range("C6").Value = [double]"9343476090.0941372"
$ws1.Range("D6").Value = [double]"21423030.825376313"
$ws1.Range("E6").Value = [double]"0.96575744716267398"
$ws1.Range("F6").Value = [double]"136.13139717353735"
$ws1.Range("G6").Value = [double]"295.18348906517315"
$ws1.Range("H6").Value = [double]"-3.6186849785866837E-2"
$ws1.Range("I6").Value = [double]"1.5122489958370178E-2"
$ws1.Range("J6").Value = [double]"5.3235774727536489E-2"
$ws1.Range("K6").Value = [double]"4.91400158657882E-2"
$ws1.Range("L6").Value = [double]"-3.0022649526117995E-4"
$ws1.Range("M6").Value = [double]"-3.087473276982311E-2"
$ws1.Range("N6").Value = [double]"101.51224899583701"
$ws1.Range("O6").Value = [double]"105.32357747275366"
$ws1.Range("A7").Value = [double]"2017"
$ws1.Range("B7").Value = [double]"452.51681152100053"
$ws1.Range("C7").Value = [double]"9204284391.8050823"
$ws1.Range("D7").Value = [double]"20340204.291786686"
$ws1.Range("E7").Value = [double]"0.92052293550703623"
$ws1.Range("F7").Value = [double]"136.17227969982335"
$ws1.Range("G7").Value = [double]"304.58754822152844"
$ws1.Range("H7").Value = [double]"-3.5869632351568481E-2"
$ws1.Range("I7").Value = [double]"8.9638823037501147E-3"
$ws1.Range("J7").Value = [double]"4.6501506600886477E-2"
$ws1.Range("K7").Value = [double]"4.2139905519349652E-2"
$ws1.Range("L7").Value = [double]"1.1321105465380787E-2"
$ws1.Range("M7").Value = [double]"-3.8898520130888015E-2"
$ws1.Range("N7").Value = [double]"100"
$ws1.Range("O7").Value = [double]"100"

# Make the data sheet the active/selected sheet and select the data rows
# (mirrors the "tabSelected" / selection state captured in the saved file).
$ws1.Activate()
$ws1.Range("A2:O7").Select()
